$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 52, shifting existing row 52 (and everything below it) down by one.
$ws.Rows(52).Insert()

# Fill the newly inserted row 52 with the new record.
$ws.Cells.Item(52, 1).Value = 4
$ws.Cells.Item(52, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(52, 3).Value = "Los Lagos"
$ws.Cells.Item(52, 4).Value = 44582
$ws.Cells.Item(52, 5).Value = 10
$ws.Cells.Item(52, 6).Value = 100112026
$ws.Cells.Item(52, 7).Value = "Haba"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 80
$ws.Cells.Item(52, 11).Value = 24000
$ws.Cells.Item(52, 12).Value = 24000
$ws.Cells.Item(52, 13).Value = 24000
$ws.Cells.Item(52, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(52, 15).Value = "Región Metropolitana"
$ws.Cells.Item(52, 16).Value = 960
$ws.Cells.Item(52, 17).Value = 25
$ws.Cells.Item(52, 18).Value = "Hortaliza"
